$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (column CK, 89th column) with the same label style as
# the rest of the header row (row 1).
$ws.Range("CK1").Value = "21-nov"

# Apply the same number formatting as the other data columns (centered,
# integer "0" format) to the new column's data rows before writing values so
# that Excel reuses the existing cell style instead of allocating a new one.
$ws.Range("CK2:CK11").HorizontalAlignment = -4108
$ws.Range("CK2:CK11").NumberFormat = "0"

# New sales figures for "21-nov" per product row.
$ws.Range("CK2").Value = 12
$ws.Range("CK3").Value = 8
$ws.Range("CK4").Value = 8
$ws.Range("CK5").Value = 13
$ws.Range("CK6").Value = 8
$ws.Range("CK7").Value = 12
$ws.Range("CK8").Value = 11
$ws.Range("CK9").Value = 12
$ws.Range("CK10").Value = 18
$ws.Range("CK11").Value = 0

# Match the author's final selection when the file was saved.
$ws.Range("CL12").Select()
